# Testdata.xlsx update: add Enums/config/Reports test-data columns, drop Sheet2.
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the now-unused "Sheet2" tab (Sheet1 and Sheet3 remain) ---
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null

$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row ---
$ws.Range("A1").Value = "testcaseName"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "expectedTitle"
$ws.Range("E1").Value = "runTestcase"

# --- tileValidationTest rows ---
$ws.Range("A2").Value = "tileValidationTest"
$ws.Range("B2").Value = "Admin"
$ws.Range("C2").Value = "admin123"
$ws.Range("D2").Value = "OrangeHRM"
$ws.Range("E2").Value = "yes"

$ws.Range("A3").Value = "tileValidationTest"
$ws.Range("B3").Value = "prabhu"
$ws.Range("C3").Value = "prabhu"
$ws.Range("D3").Value = "OrangeHRM"
$ws.Range("E3").Value = "no"

# --- Column sizing for the new/resized columns ---
# (ColumnWidth is in "characters"; the stored OOXML width is ColumnWidth + 5/6,
#  so subtract that offset to land on the desired stored width.)
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 12.592447916666666
$ws.Columns.Item(5).ColumnWidth = 10.736979166666666

# --- Selection cursor moves to E4 ---
$ws.Range("E4").Select() | Out-Null
